# Vectorsum_speed.xlsx update:
#  - Add BCC64 (Embarcadero bcc64) compiler results
#  - Rename "msvc /Od" to "msvc /Os" (and update its measurements)
#  - Update hardware/OS/compiler version info block
#  - Re-flow the options legend tables below to make room
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old layout for rows 15-60 so stale cells from the
# previous (shorter) layout don't linger after the re-flow.
$ws.Range("A15:I60").ClearContents()

# The benchmark-number columns (B:I) for the compiler rows use the
# thousands-separator number format; make sure the newly added bcc64
# rows (20-24) pick up the same formatting as the existing rows.
$ws.Range("B15:I24").NumberFormat = "#,##0"

$ws.Range("A20").Value = "bcc64 -Og"
$ws.Range("A21").Value = "bcc64 -O1"
$ws.Range("A22").Value = "bcc64 -O2"
$ws.Range("A23").Value = "bcc64 -O3"
$ws.Range("A24").Value = "bcc64 -O3 -march=haswell"
$ws.Range("A15").Value = "msvc /Os"
$ws.Range("B30").Value = "Windows 11 Home"
$ws.Range("D30").Value = "10.0.22000.613"
$ws.Range("A30").Value = "Visual Studio 2022,  19.31.31107"
$ws.Range("A31").Value = "Embarcadero C++ 7.50 for Win64"
$ws.Range("B31").Value = " bcc64 version 5.0.2 (d939c99b.1e953156.37585) (based on LLVM 5.0.2)"
$ws.Range("A60").Value = "Embarcadero bcc64 options like Clang"
$ws.Range("B15").Value = 3268
$ws.Range("C15").Value = 3915
$ws.Range("D15").Value = 10708
$ws.Range("E15").Value = 1704
$ws.Range("F15").Value = 1482
$ws.Range("G15").Value = 1598
$ws.Range("H15").Value = 2854
$ws.Range("I15").Value = 1980
$ws.Range("A16").Value = "msvc /O1"
$ws.Range("B16").Value = 506
$ws.Range("C16").Value = 363
$ws.Range("D16").Value = 348
$ws.Range("E16").Value = 358
$ws.Range("F16").Value = 358
$ws.Range("G16").Value = 514
$ws.Range("H16").Value = 364
$ws.Range("I16").Value = 366
$ws.Range("A17").Value = "msvc /O2"
$ws.Range("B17").Value = 165
$ws.Range("C17").Value = 365
$ws.Range("D17").Value = 349
$ws.Range("E17").Value = 169
$ws.Range("F17").Value = 164
$ws.Range("G17").Value = 164
$ws.Range("H17").Value = 168
$ws.Range("I17").Value = 364
$ws.Range("A18").Value = "msvc /O2 /Ob3"
$ws.Range("B18").Value = 163
$ws.Range("C18").Value = 358
$ws.Range("D18").Value = 347
$ws.Range("E18").Value = 165
$ws.Range("F18").Value = 165
$ws.Range("G18").Value = 166
$ws.Range("H18").Value = 166
$ws.Range("I18").Value = 359
$ws.Range("A19").Value = "msvc /O2 /Ob3 /arch:avx2"
$ws.Range("B19").Value = 148
$ws.Range("C19").Value = 359
$ws.Range("D19").Value = 347
$ws.Range("E19").Value = 149
$ws.Range("F19").Value = 149
$ws.Range("G19").Value = 154
$ws.Range("H19").Value = 154
$ws.Range("I19").Value = 362
$ws.Range("B20").Value = 9266
$ws.Range("C20").Value = 14486
$ws.Range("D20").Value = 4739
$ws.Range("E20").Value = 6662
$ws.Range("F20").Value = 361
$ws.Range("G20").Value = 369
$ws.Range("H20").Value = 1415
$ws.Range("B21").Value = 9234
$ws.Range("C21").Value = 14310
$ws.Range("D21").Value = 4724
$ws.Range("E21").Value = 6598
$ws.Range("F21").Value = 357
$ws.Range("G21").Value = 370
$ws.Range("H21").Value = 1369
$ws.Range("B22").Value = 161
$ws.Range("C22").Value = 160
$ws.Range("D22").Value = 161
$ws.Range("E22").Value = 161
$ws.Range("F22").Value = 161
$ws.Range("G22").Value = 161
$ws.Range("H22").Value = 161
$ws.Range("B23").Value = 163
$ws.Range("C23").Value = 165
$ws.Range("D23").Value = 164
$ws.Range("E23").Value = 162
$ws.Range("F23").Value = 161
$ws.Range("G23").Value = 162
$ws.Range("H23").Value = 162
$ws.Range("B24").Value = 150
$ws.Range("C24").Value = 152
$ws.Range("D24").Value = 154
$ws.Range("E24").Value = 151
$ws.Range("F24").Value = 153
$ws.Range("G24").Value = 150
$ws.Range("H24").Value = 150
$ws.Range("A28").Value = "gcc 11.2"
$ws.Range("B28").Value = "Linux Fedora 35"
$ws.Range("A29").Value = "clang 13.0"
$ws.Range("B29").Value = "Linux Fedora 35"
$ws.Range("A34").Value = "Hardware"
$ws.Range("A35").Value = "CPU Intel Core i7-8700, 3.2 GHz"
$ws.Range("B35").Value = "6 cores, 12 threads"
$ws.Range("D35").Value = "8th gen, Coffee Lake"
$ws.Range("A36").Value = "32 GB RAM"
$ws.Range("A39").Value = "MSVC options:"
$ws.Range("A40").Value = " /Od"
$ws.Range("B40").Value = "Disables optimization"
$ws.Range("A41").Value = " /O1"
$ws.Range("B41").Value = "Creates small code"
$ws.Range("A42").Value = " /O2"
$ws.Range("B42").Value = "Creates fast code"
$ws.Range("A43").Value = " /Ob3"
$ws.Range("B43").Value = "Aggressive inlining"
$ws.Range("A44").Value = " /arch:avx2"
$ws.Range("B44").Value = "Enables the use of Intel Advanced Vector Extensions 2 instructions."
$ws.Range("A46").Value = "GCC options"
$ws.Range("A47").Value = " -Og"
$ws.Range("B47").Value = "Optimize debugging experience"
$ws.Range("A48").Value = " -O1"
$ws.Range("B48").Value = "Optimize"
$ws.Range("A49").Value = " -O2"
$ws.Range("B49").Value = "Optimize even more, nearly all optimizations that do not involve a space-speed tradeoff"
$ws.Range("A50").Value = " -O3"
$ws.Range("B50").Value = "Optimize yet more"
$ws.Range("A51").Value = " -march=haswell"
$ws.Range("B51").Value = "Enable instructions found on Haswell and later CPUs"
$ws.Range("A53").Value = "Clang options"
$ws.Range("A54").Value = " -Og"
$ws.Range("B54").Value = "Like -O1"
$ws.Range("A55").Value = " -O1"
$ws.Range("B55").Value = "Optimize"
$ws.Range("A56").Value = " -O2"
$ws.Range("B56").Value = "Moderate level of optimization which enables most optimizations"
$ws.Range("A57").Value = " -O3"
$ws.Range("B57").Value = "Like -O2, except that it enables optimizations that take longer to perform or that may generate larger code (in an attempt to make the program run faster)"
$ws.Range("A58").Value = " -march=haswell"
$ws.Range("B58").Value = "Enable instructions found on Haswell and later CPUs"

# Update the selected cell to match the saved worksheet state.
$ws.Range("H28").Select()
